$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Comments" column header (D1) to the new header text.
$ws.Range("D1").Value = "Abou the Agent Learning"

# Update the view: scroll back to top-left and move the active selection.
$ws.Range("I8").Select()
